$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-05 19:17:49"
$ws.Range("H2").Value = "'95%"
$ws.Range("E3").Value = "2026-02-05 19:17:52"
$ws.Range("E4").Value = "2026-02-05 19:17:54"
$ws.Range("J4").Value = "989.9 hPa"
$ws.Range("O4").Value = "11.3 °C"
$ws.Range("E5").Value = "2026-02-05 19:17:57"
$ws.Range("J5").Value = "990.1 hPa"
$ws.Range("O5").Value = "9.7 °C"
$ws.Range("E6").Value = "2026-02-05 19:17:59"
$ws.Range("H6").Value = "'72%"
$ws.Range("J6").Value = "991.8 hPa"
$ws.Range("M6").Value = "16.6 °C 18:33 TU"
$ws.Range("O6").Value = "12.8 °C"
$ws.Range("E7").Value = "2026-02-05 19:18:02"
$ws.Range("E8").Value = "2026-02-05 19:18:04"
$ws.Range("H8").Value = "'87%"
$ws.Range("M8").Value = "14.7 °C 18:31 TU"
$ws.Range("O8").Value = "8.6 °C"
$ws.Range("E9").Value = "2026-02-05 19:18:07"
$ws.Range("E10").Value = "2026-02-05 19:18:09"
$ws.Range("E11").Value = "2026-02-05 19:18:11"
$ws.Range("J11").Value = "994.9 hPa"
$ws.Range("O11").Value = "0.4 °C"
$ws.Range("E12").Value = "2026-02-05 19:18:14"
$ws.Range("H12").Value = "'89%"
$ws.Range("O12").Value = "9.9 °C"
$ws.Range("E13").Value = "2026-02-05 19:18:17"
$ws.Range("E14").Value = "2026-02-05 19:18:19"
$ws.Range("I14").Value = "6.8 mm"
$ws.Range("O14").Value = "-2.3 °C"
$ws.Range("E15").Value = "2026-02-05 19:18:21"
$ws.Range("H15").Value = "'83%"
$ws.Range("J15").Value = "990.6 hPa"
$ws.Range("O15").Value = "8.1 °C"
$ws.Range("E16").Value = "2026-02-05 19:18:24"
$ws.Range("O16").Value = "3.7 °C"
$ws.Range("E17").Value = "2026-02-05 19:18:27"
$ws.Range("I17").Value = "8.5 mm"
$ws.Range("J17").Value = "995.2 hPa"
$ws.Range("O17").Value = "0.9 °C"
$ws.Range("E18").Value = "2026-02-05 19:18:29"
$ws.Range("E19").Value = "2026-02-05 19:18:32"
$ws.Range("O19").Value = "7.2 °C"
$ws.Range("E20").Value = "2026-02-05 19:18:34"
$ws.Range("K20").Value = "1.2 MJ/m2"
$ws.Range("E21").Value = "2026-02-05 19:18:37"
$ws.Range("J21").Value = "990.8 hPa"
$ws.Range("O21").Value = "6.2 °C"
$ws.Range("E22").Value = "2026-02-05 19:18:39"
$ws.Range("H22").Value = "'89%"
$ws.Range("M22").Value = "16.0 °C 18:50 TU"
$ws.Range("O22").Value = "8.5 °C"
$ws.Range("E23").Value = "2026-02-05 19:18:42"
$ws.Range("J23").Value = "990.0 hPa"
$ws.Range("E24").Value = "2026-02-05 19:18:44"
$ws.Range("J24").Value = "989.1 hPa"
$ws.Range("E25").Value = "2026-02-05 19:18:47"
$ws.Range("J25").Value = "994.1 hPa"
$ws.Range("M25").Value = "3.2 °C 18:59 TU"
$ws.Range("O25").Value = "0.6 °C"
$ws.Range("E26").Value = "2026-02-05 19:18:49"
$ws.Range("H26").Value = "'76%"
$ws.Range("E27").Value = "2026-02-05 19:18:52"
$ws.Range("J27").Value = "990.3 hPa"
$ws.Range("E28").Value = "2026-02-05 19:18:54"
$ws.Range("H28").Value = "'95%"
$ws.Range("J28").Value = "993.1 hPa"
$ws.Range("O28").Value = "2.4 °C"
$ws.Range("E29").Value = "2026-02-05 19:18:57"
$ws.Range("H29").Value = "'81%"
$ws.Range("O29").Value = "9.0 °C"
$ws.Range("E30").Value = "2026-02-05 19:18:59"
$ws.Range("I30").Value = "5.1 mm"
$ws.Range("O30").Value = "-1.8 °C"
$ws.Range("E31").Value = "2026-02-05 19:19:02"
$ws.Range("I31").Value = "18.8 mm"
$ws.Range("J31").Value = "994.4 hPa"
$ws.Range("E32").Value = "2026-02-05 19:19:04"
$ws.Range("E33").Value = "2026-02-05 19:19:07"
$ws.Range("H33").Value = "'87%"
$ws.Range("O33").Value = "9.0 °C"
$ws.Range("E34").Value = "2026-02-05 19:19:09"
$ws.Range("O34").Value = "3.8 °C"
$ws.Range("E35").Value = "2026-02-05 19:19:12"
$ws.Range("I35").Value = "4.8 mm"
$ws.Range("E36").Value = "2026-02-05 19:19:14"
